$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "image size" inputs that now drive the chart centre (B1/C1).
$ws.Range("I1").Value = "Image Size:"
$ws.Range("J1").Value = 1961
$ws.Range("J2").Value = 1787

# Centre (B1, C1) is now derived from the image size instead of being a
# hard-coded literal.
$ws.Range("B1").Formula = "=J1/2"
$ws.Range("C1").Formula = "=J2/2"

# Radius grew considerably, which is what actually fixes the clickable
# <area> coordinates for the radar chart image map.
$ws.Range("B2").Value = 1200

# Selection moved off the H5:H11 helper column onto B2 and the sheet no
# longer scrolls so B1 is the left-most visible column.
$ws.Range("B2").Select()
